$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Server Info = server version + OS information " -> split into several
#    runs ("Server Info = " / "Id + " / object_vendor_id / " + " /
#    object_product_id / " + " / object_version_id / " "), with
#    w:proofErr spellStart/spellEnd markers wrapping each object_* run
#    (Word's automatic "flag possible misspelling" markup).
# ---------------------------------------------------------------------------

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Server Info = server version + OS information `r") {
        $target = $p
    }
}

$full = $target.Range
# Range up to (but excluding) the trailing paragraph mark, so InsertXML only
# rewrites the runs inside the paragraph and leaves the <w:p>/<w:pPr> (and
# their rsid/style attributes) untouched.
$r = $d.Range($full.Start, $full.End - 1)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Server Info = </w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Id + </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>object_vendor_id</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> + </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>object_product_id</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> + </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>object_version_id</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) Remove the whole "Issues" paragraph that used to follow the
#    "Unresolved Issue" heading paragraph.
# ---------------------------------------------------------------------------

$issues = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Issues`r") {
        $issues = $p
    }
}

$issues.Range.Delete()
